$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.257751226425171
$ws.Range("B1").Value = 4.961606502532959
$ws.Range("C1").Value = 2.104001998901367
$ws.Range("D1").Value = 1.513599157333374
$ws.Range("E1").Value = 1.301970720291138
